# Apply updated cryptocurrency price/volume data to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B/C (name/link) and E (volume %) cells are never numeric-looking,
# so a plain Value assignment keeps them as text.
$textUpdates = @(
    @{Cell = "E2"; Value = '  -1.71%  '}
    @{Cell = "E3"; Value = '  -0.83%  '}
    @{Cell = "E4"; Value = '  -0.05%  '}
    @{Cell = "E5"; Value = '  -0.87%  '}
    @{Cell = "E6"; Value = '  -2.26%  '}
    @{Cell = "E7"; Value = '  -0.13%  '}
    @{Cell = "E8"; Value = '  +0.21%  '}
    @{Cell = "E9"; Value = '  -0.93%  '}
    @{Cell = "E10"; Value = '  -1.27%  '}
    @{Cell = "E11"; Value = '  -0.58%  '}
    @{Cell = "E12"; Value = '  +3.35%  '}
    @{Cell = "E13"; Value = '  +1.19%  '}
    @{Cell = "E14"; Value = '  -0.88%  '}
    @{Cell = "E15"; Value = '  -1.84%  '}
    @{Cell = "E16"; Value = '  -2.78%  '}
    @{Cell = "E17"; Value = '  -1.51%  '}
    @{Cell = "E18"; Value = '  -0.82%  '}
    @{Cell = "E19"; Value = '  -2.19%  '}
    @{Cell = "E20"; Value = '  -0.84%  '}
    @{Cell = "E21"; Value = '  -2.24%  '}
    @{Cell = "E22"; Value = '  -1.52%  '}
    @{Cell = "E23"; Value = '  +0.44%  '}
    @{Cell = "E24"; Value = '  -0.05%  '}
    @{Cell = "E25"; Value = '  +1.10%  '}
    @{Cell = "E26"; Value = '  -4.20%  '}
    @{Cell = "E27"; Value = '  -3.83%  '}
    @{Cell = "E28"; Value = '  -0.07%  '}
    @{Cell = "E29"; Value = '  -0.63%  '}
    @{Cell = "E30"; Value = '  -3.46%  '}
    @{Cell = "E31"; Value = '  -1.16%  '}
    @{Cell = "E32"; Value = '  -5.16%  '}
    @{Cell = "E33"; Value = '  -6.78%  '}
    @{Cell = "E34"; Value = '  -2.49%  '}
    @{Cell = "E35"; Value = '  -0.12%  '}
    @{Cell = "B36"; Value = 'Kaspa'}
    @{Cell = "C36"; Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'}
    @{Cell = "E36"; Value = '  -0.60%  '}
    @{Cell = "B37"; Value = 'Monero'}
    @{Cell = "C37"; Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'}
    @{Cell = "E37"; Value = '  +0.07%  '}
    @{Cell = "E38"; Value = '  +0.69%  '}
    @{Cell = "E39"; Value = '  -1.04%  '}
    @{Cell = "E40"; Value = '  -0.06%  '}
    @{Cell = "E41"; Value = '  -1.33%  '}
    @{Cell = "E42"; Value = '  -1.42%  '}
    @{Cell = "E43"; Value = '  -3.15%  '}
    @{Cell = "E44"; Value = '  -0.85%  '}
    @{Cell = "E45"; Value = '  -8.09%  '}
    @{Cell = "E46"; Value = '  -7.83%  '}
    @{Cell = "E47"; Value = '  -1.07%  '}
    @{Cell = "E48"; Value = '  -2.11%  '}
    @{Cell = "E49"; Value = '  -2.12%  '}
    @{Cell = "E50"; Value = '  -0.07%  '}
    @{Cell = "E51"; Value = '  -4.00%  '}
)
foreach ($u in $textUpdates) {
    $ws.Range($u.Cell).Value = $u.Value
}

# Column D (price) values often look numeric (e.g. "0.999"), so Excel would
# auto-convert them to numbers. Force text storage via NumberFormat, then
# restore the default "Normal" style so no stray formatting is left behind.
$priceUpdates = @(
    @{Cell = "D2"; Value = '69.273.77'}
    @{Cell = "D3"; Value = '2.504.61'}
    @{Cell = "D4"; Value = '0.999'}
    @{Cell = "D5"; Value = '571.57'}
    @{Cell = "D6"; Value = '165.47'}
    @{Cell = "D8"; Value = '0.512'}
    @{Cell = "D9"; Value = '2.502.65'}
    @{Cell = "D10"; Value = '0.158'}
    @{Cell = "D12"; Value = '0.355'}
    @{Cell = "D13"; Value = '4.89'}
    @{Cell = "D14"; Value = '2.962.38'}
    @{Cell = "D15"; Value = '69.023.31'}
    @{Cell = "D17"; Value = '24.71'}
    @{Cell = "D18"; Value = '2.494.24'}
    @{Cell = "D19"; Value = '11.26'}
    @{Cell = "D20"; Value = '7.57'}
    @{Cell = "D21"; Value = '347.79'}
    @{Cell = "D22"; Value = '3.90'}
    @{Cell = "D23"; Value = '1.99'}
    @{Cell = "D25"; Value = '70.06'}
    @{Cell = "D26"; Value = '3.92'}
    @{Cell = "D27"; Value = '8.84'}
    @{Cell = "D28"; Value = '2.649.69'}
    @{Cell = "D30"; Value = '0.0₃0884'}
    @{Cell = "D31"; Value = '7.78'}
    @{Cell = "D32"; Value = '459.04'}
    @{Cell = "D33"; Value = '1.23'}
    @{Cell = "D35"; Value = '0.998'}
    @{Cell = "D36"; Value = '0.116'}
    @{Cell = "D37"; Value = '157.19'}
    @{Cell = "D38"; Value = '19.03'}
    @{Cell = "D39"; Value = '18.44'}
    @{Cell = "D41"; Value = '0.317'}
    @{Cell = "D42"; Value = '4.68'}
    @{Cell = "D43"; Value = '1.60'}
    @{Cell = "D44"; Value = '38.01'}
    @{Cell = "D46"; Value = '2.20'}
    @{Cell = "D47"; Value = '141.65'}
    @{Cell = "D48"; Value = '3.46'}
    @{Cell = "D49"; Value = '0.517'}
    @{Cell = "D50"; Value = '0.0729'}
    @{Cell = "D51"; Value = '0.576'}
)
foreach ($u in $priceUpdates) {
    $cellRange = $ws.Range($u.Cell)
    $cellRange.NumberFormat = "@"
    $cellRange.Value = $u.Value
    $cellRange.Style = "Normal"
}
